# Adds the "For The Philanthropist" portfolio item after the AirHop
# item, and relocates the stray "_GoBack" bookmark (which had been
# sitting mid-sentence in the AirHop paragraph) to the end of the
# newly-typed content -- mirroring where Word leaves it after the
# last edit.

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Two blank paragraphs directly after the (pre-existing) AirHop
#    paragraph. We never touch Font on a collapsed range that still
#    belongs to the *original* paragraph, so this first
#    InsertParagraphAfter() call is done with no formatting changes.
# -----------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# -----------------------------------------------------------------
# 2. Bold heading paragraph: "For The Philanthropist Site"
#    (format AFTER inserting each chunk of text).
# -----------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertAfter("For The Philanthropist")
$r.Font.Bold = $true

$r.Collapse(0)
$r.InsertAfter(" Site")
$r.Font.Bold = $true

# -----------------------------------------------------------------
# 3. Blank paragraph.
# -----------------------------------------------------------------
$r.Collapse(0)
$r.InsertParagraphAfter()

# -----------------------------------------------------------------
# 4. Body paragraph introducing For The Philanthropist.
# -----------------------------------------------------------------
$paraC = @(
    "For The Philanthropist",
    ", ",
    "a company devoted to ",
    "engaging ",
    "and assisting non-profits ",
    "in ",
    "meeting ",
    "their ",
    "philanthropic goals",
    " ",
    "and ",
    "enacting social change, ",
    "approached ",
    "me at the beginning of the summer of ",
    "2014",
    " to ask ",
    "me to develop ",
    "its",
    " new website.",
    " ",
    "The company",
    " primarily ",
    "wanted an informational ",
    "site",
    " that advertised their work in a clean",
    " and ",
    "intuitive way."
)

$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
foreach ($chunk in $paraC) {
    $r.InsertAfter($chunk)
    $r.Font.Bold = $false
    $r.Collapse(0)
}

# -----------------------------------------------------------------
# 5. Blank paragraph.
# -----------------------------------------------------------------
$r.InsertParagraphAfter()

# -----------------------------------------------------------------
# 6. Closing paragraph about the site design, up through
#    "minimalistic" -- this is where the author's cursor was left,
#    so the "_GoBack" bookmark is re-anchored there before the final
#    " " + "themes." runs are typed.
# -----------------------------------------------------------------
$paraD1 = @(
    "To ",
    "complement For The Philanthropist’s",
    " noble and ",
    "morally clean",
    " business",
    " goals",
    ", the interfaced",
    " I designed ",
    "features ",
    "a simple",
    " structure ",
    "with carefully ",
    "picked",
    " font faces",
    " and ",
    "minimalistic"
)

$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
foreach ($chunk in $paraD1) {
    $r.InsertAfter($chunk)
    $r.Font.Bold = $false
    $r.Collapse(0)
}

# Relocate "_GoBack" here -- right before the final " " + "themes."
# runs. Bookmark names are unique, so re-adding it removes it from
# its old spot (mid-sentence in the AirHop paragraph).
$d.Bookmarks.Add("_GoBack", $r) | Out-Null

$r.InsertAfter(" ")
$r.Font.Bold = $false
$r.Collapse(0)

$r.InsertAfter("themes.")
$r.Font.Bold = $false
$r.Collapse(0)

Write-Output "done"
